$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1780495
$ws.Range("H2").Value = 0.356099
$ws.Range("M2").Value = 8.949653
$ws.Range("N2").Value = 17.899306
$ws.Range("O2").Value = 0.1668927877080592
$ws.Range("P2").Value = 0.1610232428880788
$ws.Range("Q2").Value = 1.5934812418235
$ws.Range("R2").Value = 6.373924967294
$ws.Range("S2").Value = 0.1668927877080592
$ws.Range("T2").Value = 0.1610232428880788

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1780495
$ws.Range("H3").Value = 0.356099
$ws.Range("O3").Value = 0.06249788578732534
$ws.Range("P3").Value = 0.09044979457765322
$ws.Range("Q3").Value = 0.5967256585703333
$ws.Range("R3").Value = 3.580353951422
$ws.Range("S3").Value = 0.06249788578732534
$ws.Range("T3").Value = 0.09044979457765322

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.1780495
$ws.Range("H4").Value = 0.356099
$ws.Range("M4").Value = 0.050758
$ws.Range("N4").Value = 0.152274
$ws.Range("O4").Value = 0.0009465332475444208
$ws.Range("P4").Value = 0.001369866143834812
$ws.Range("Q4").Value = 0.009037436520999999
$ws.Range("R4").Value = 0.054224619126
$ws.Range("S4").Value = 0.0009465332475444208
$ws.Range("T4").Value = 0.001369866143834812

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.1780495
$ws.Range("H5").Value = 0.356099
$ws.Range("M5").Value = 40.766071
$ws.Range("N5").Value = 81.53214199999999
$ws.Range("O5").Value = 0.7602041367519689
$ws.Range("P5").Value = 0.7334680967212543
$ws.Range("Q5").Value = 7.258378558514499
$ws.Range("R5").Value = 29.033514234058
$ws.Range("S5").Value = 0.7602041367519689
$ws.Range("T5").Value = 0.7334680967212543

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.1780495
$ws.Range("H6").Value = 0.356099
$ws.Range("M6").Value = 0.3317233333333334
$ws.Range("N6").Value = 0.9951700000000001
$ws.Range("O6").Value = 0.006185964064507279
$ws.Range("P6").Value = 0.008952609705925435
$ws.Range("Q6").Value = 0.05906317363833334
$ws.Range("R6").Value = 0.35437904183
$ws.Range("S6").Value = 0.006185964064507279
$ws.Range("T6").Value = 0.008952609705925435

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.1780495
$ws.Range("H7").Value = 0.356099
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1754986666666667
$ws.Range("N7").Value = 0.526496
$ws.Range("O7").Value = 0.003272692440594897
$ws.Range("P7").Value = 0.004736389963253432
$ws.Range("Q7").Value = 0.03124744985066667
$ws.Range("R7").Value = 0.187484699104
$ws.Range("S7").Value = 0.003272692440594897
$ws.Range("T7").Value = 0.004736389963253432

Write-Host "Applied TPM updates"
